# 自动更新Excel文件 - daily decrement of remaining-days column (E),
# with auto-refill of the cycle (E reset to total days D, start date F
# reset to today) whenever the countdown would reach zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStartDate = 20251107

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $total = $dCell.Value()
    $remaining = $eCell.Value()
    $startDate = $fCell.Value()

    if ($remaining -eq $null -or $total -eq $null) {
        continue
    }

    # Skip rows whose start date is not a well-formed yyyymmdd value
    # (corrupt source data - left untouched by the real update job).
    $dateText = [string]([int]$startDate)
    if ($dateText.Length -ne 8) {
        continue
    }

    $newRemaining = $remaining - 1

    if ($newRemaining -le 0) {
        $eCell.Value = $total
        $fCell.Value = $newStartDate
    } else {
        $eCell.Value = $newRemaining
    }
}
